$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric (single-decimal strings like "318.51")
# would otherwise be auto-coerced to a Double by Excel on assignment, losing
# their original text-string cell type. Force them to remain text first.
$textForceCells = @(
    @{ Ref = 'D5'; Value = '318.51' }
    @{ Ref = 'D6'; Value = '1.003' }
    @{ Ref = 'D7'; Value = '0.4333' }
    @{ Ref = 'D8'; Value = '0.3730' }
    @{ Ref = 'D9'; Value = '0.07428' }
    @{ Ref = 'D10'; Value = '0.9307' }
    @{ Ref = 'D11'; Value = '21.21' }
    @{ Ref = 'D14'; Value = '5.429' }
    @{ Ref = 'D15'; Value = '0.06856' }
    @{ Ref = 'D16'; Value = '1.005' }
    @{ Ref = 'D17'; Value = '80.75' }
    @{ Ref = 'D18'; Value = '0.000009023' }
    @{ Ref = 'D19'; Value = '1.003' }
    @{ Ref = 'D22'; Value = '5.122' }
    @{ Ref = 'D23'; Value = '11.00' }
    @{ Ref = 'D25'; Value = '2.041' }
    @{ Ref = 'D26'; Value = '153.80' }
    @{ Ref = 'D27'; Value = '18.50' }
    @{ Ref = 'D28'; Value = '5.506' }
    @{ Ref = 'D29'; Value = '112.97' }
    @{ Ref = 'D30'; Value = '1.694' }
    @{ Ref = 'D31'; Value = '0.08980' }
    @{ Ref = 'D32'; Value = '0.8073' }
    @{ Ref = 'D33'; Value = '4.781' }
    @{ Ref = 'D34'; Value = '1.175' }
    @{ Ref = 'D35'; Value = '2.949' }
    @{ Ref = 'D38'; Value = '0.05486' }
    @{ Ref = 'D39'; Value = '0.01969' }
    @{ Ref = 'D40'; Value = '2.991' }
    @{ Ref = 'D41'; Value = '0.5236' }
    @{ Ref = 'D42'; Value = '7.002' }
    @{ Ref = 'D43'; Value = '0.1683' }
    @{ Ref = 'D44'; Value = '8.759' }
    @{ Ref = 'D45'; Value = '0.06720' }
    @{ Ref = 'D46'; Value = '0.4873' }
    @{ Ref = 'D47'; Value = '10.55' }
    @{ Ref = 'D48'; Value = '106.79' }
    @{ Ref = 'D50'; Value = '1.672' }
    @{ Ref = 'D51'; Value = '1.860' }
)
foreach ($item in $textForceCells) {
    $rng = $ws.Range($item.Ref)
    $rng.NumberFormat = "@"
    $rng.Value = $item.Value
}

# Cells whose new values are already safely non-numeric text (e.g. contain
# two separator dots like "27.988.57", or padded percentages) can be set directly.
$directCells = @(
    @{ Ref = 'D2'; Value = '27.988.57' }
    @{ Ref = 'E2'; Value = '  -3.80%  ' }
    @{ Ref = 'D3'; Value = '1.867.83' }
    @{ Ref = 'E3'; Value = '  -2.89%  ' }
    @{ Ref = 'E4'; Value = '  -0.25%  ' }
    @{ Ref = 'E5'; Value = '  -2.35%  ' }
    @{ Ref = 'E7'; Value = '  -5.97%  ' }
    @{ Ref = 'E8'; Value = '  -2.44%  ' }
    @{ Ref = 'E9'; Value = '  -4.38%  ' }
    @{ Ref = 'E10'; Value = '  -4.90%  ' }
    @{ Ref = 'E11'; Value = '  -6.53%  ' }
    @{ Ref = 'D12'; Value = '1.981.08' }
    @{ Ref = 'E12'; Value = '  +0.88%  ' }
    @{ Ref = 'E13'; Value = '  -3.41%  ' }
    @{ Ref = 'E14'; Value = '  -4.67%  ' }
    @{ Ref = 'E15'; Value = '  -2.58%  ' }
    @{ Ref = 'E16'; Value = '  -0.35%  ' }
    @{ Ref = 'E17'; Value = '  -4.28%  ' }
    @{ Ref = 'E18'; Value = '  -4.94%  ' }
    @{ Ref = 'E19'; Value = '  -0.34%  ' }
    @{ Ref = 'E20'; Value = '  -6.03%  ' }
    @{ Ref = 'D21'; Value = '27.982.70' }
    @{ Ref = 'E21'; Value = '  -3.88%  ' }
    @{ Ref = 'E22'; Value = '  -4.13%  ' }
    @{ Ref = 'E23'; Value = '  +0.28%  ' }
    @{ Ref = 'D24'; Value = '2.179.81' }
    @{ Ref = 'E24'; Value = '  +1.90%  ' }
    @{ Ref = 'E25'; Value = '  -1.37%  ' }
    @{ Ref = 'E26'; Value = '  -2.73%  ' }
    @{ Ref = 'E27'; Value = '  -2.89%  ' }
    @{ Ref = 'E28'; Value = '  -2.67%  ' }
    @{ Ref = 'E29'; Value = '  -4.17%  ' }
    @{ Ref = 'E30'; Value = '  -7.71%  ' }
    @{ Ref = 'E31'; Value = '  -3.92%  ' }
    @{ Ref = 'E32'; Value = '  -5.57%  ' }
    @{ Ref = 'E33'; Value = '  -6.51%  ' }
    @{ Ref = 'E34'; Value = '  -5.39%  ' }
    @{ Ref = 'E35'; Value = '  -2.47%  ' }
    @{ Ref = 'E36'; Value = '  -0.26%  ' }
    @{ Ref = 'E37'; Value = '  -3.48%  ' }
    @{ Ref = 'E38'; Value = '  -3.47%  ' }
    @{ Ref = 'E39'; Value = '  -3.74%  ' }
    @{ Ref = 'E40'; Value = '  -6.14%  ' }
    @{ Ref = 'E41'; Value = '  -5.04%  ' }
    @{ Ref = 'E42'; Value = '  -6.64%  ' }
    @{ Ref = 'E43'; Value = '  -4.21%  ' }
    @{ Ref = 'E44'; Value = '  -6.09%  ' }
    @{ Ref = 'E45'; Value = '  -3.02%  ' }
    @{ Ref = 'E46'; Value = '  -6.42%  ' }
    @{ Ref = 'E47'; Value = '  -6.32%  ' }
    @{ Ref = 'E48'; Value = '  -3.29%  ' }
    @{ Ref = 'E49'; Value = '  -0.34%  ' }
    @{ Ref = 'E50'; Value = '  -5.36%  ' }
    @{ Ref = 'E51'; Value = '  -15.12%  ' }
)
foreach ($item in $directCells) {
    $ws.Range($item.Ref).Value = $item.Value
}
